$d = $word.ActiveDocument

# The document contains citation lines that read "Source: <url>" as a
# single run. The author wants these re-labelled "Project: <url>" so the
# wording matches the other citation lines in the document, but the run
# needs to end up split into two runs -- "Project" and ": " -- matching
# the canonical XML (the colon/space keeps its own run rather than being
# merged into the label run).

while ($true) {
    $rng = $d.Content
    $found = $rng.Find.Execute("Source: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }

    $start = $rng.Start

    # Replace "Source: " with "Project: " as a single run first.
    $rng.Text = "Project: "

    # Now split "Project: " into two runs -- "Project" and ": " -- by
    # nudging the formatting of the "Project" portion away and back,
    # which forces the engine to materialize it as its own run even
    # though the effective formatting is unchanged.
    $label = $d.Range($start, $start + 7)
    $origSize = $label.Font.Size
    $label.Font.Size = $origSize + 1
    $label.Font.Size = $origSize
}
